$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.021.31'
$ws.Range("E2").Value = '  -0.34%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.642.93'
$ws.Range("E3").Value = '  -1.38%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '207.08'
$ws.Range("E5").Value = '  -1.26%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5153'
$ws.Range("E6").Value = '  -0.51%  '
$ws.Range("E7").Value = '  +0.18%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2565'
$ws.Range("E8").Value = '  -2.45%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06204'
$ws.Range("E9").Value = '  -0.23%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.54'
$ws.Range("E10").Value = '  -2.60%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07568'
$ws.Range("E11").Value = '  +1.09%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.646.05'
$ws.Range("E12").Value = '  -1.24%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.362'
$ws.Range("E13").Value = '  -1.16%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.867.29'
$ws.Range("E14").Value = '  -1.34%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5344'
$ws.Range("E15").Value = '  -4.18%  '
$ws.Range("E16").Value = '  +1.40%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '65.64'
$ws.Range("E17").Value = '  -0.41%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '26.026.26'
$ws.Range("E18").Value = '  -0.43%  '
$ws.Range("E19").Value = '  +0.07%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.643'
$ws.Range("E20").Value = '  -2.66%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '186.04'
$ws.Range("E21").Value = '  -0.07%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.01'
$ws.Range("E22").Value = '  -3.29%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.093'
$ws.Range("E23").Value = '  -1.15%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.005'
$ws.Range("E24").Value = '  +0.19%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '147.41'
$ws.Range("E25").Value = '  -0.17%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1197'
$ws.Range("E26").Value = '  -3.52%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.339'
$ws.Range("E27").Value = '  -2.64%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.51'
$ws.Range("E28").Value = '  -2.47%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.357'
$ws.Range("E29").Value = '  -0.14%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.06010'
$ws.Range("E30").Value = '  -3.63%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.244'
$ws.Range("E31").Value = '  -2.10%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.410'
$ws.Range("E32").Value = '  -1.67%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.380'
$ws.Range("E33").Value = '  -0.99%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.615'
$ws.Range("E34").Value = '  -0.23%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9683'
$ws.Range("E35").Value = '  -2.65%  '
$ws.Range("E36").Value = '  -0.91%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.731'
$ws.Range("E37").Value = '  +1.10%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5840'
$ws.Range("E38").Value = '  -2.89%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01589'
$ws.Range("E39").Value = '  -0.77%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.073.72'
$ws.Range("E40").Value = '  +0.08%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.831'
$ws.Range("E41").Value = '  -4.65%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8439'
$ws.Range("E42").Value = '  -1.89%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.003'
$ws.Range("E43").Value = '  +0.02%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '100.23'
$ws.Range("E44").Value = '  +1.22%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.797.78'
$ws.Range("E45").Value = '  -0.86%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0₈109'
$ws.Range("E46").Value = '  -0.55%  '
$ws.Range("E47").Value = '  -0.55%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '54.37'
$ws.Range("E48").Value = '  -2.70%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.940'
$ws.Range("E49").Value = '  +0.46%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05209'
$ws.Range("E50").Value = '  -0.80%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4236'
$ws.Range("E51").Value = '  -0.22%  '

Write-Output "Updated cryptos list"
